$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as literal text in this workbook (e.g.
# "317.30", "45.350.35"), not as numbers. Excel auto-converts a plain numeric-
# looking string assigned via .Value into a real Number, which would both change
# the cell type and silently drop significant trailing zeros (e.g. "0.0920" ->
# 0.092). To keep every updated Price cell a literal text value matching the
# source data, force Text format before writing, then restore the default
# (Normal) style afterwards so no stray number-format is left behind.
$textForceRefs = @('D2', 'D3', 'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '45.468.65'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.381.15'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '319.28'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '109.68'
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('D7').Value = '0.637'
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').Value = '41.07'
$ws.Range('E10').Value = '  -3.74%  '
$ws.Range('D11').Value = '0.0920'
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('D12').Value = '8.57'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '0.988'
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('D15').Value = '15.60'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '2.741.92'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '2.377.56'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '45.376.57'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = '15.71'
$ws.Range('E19').Value = '  +16.91%  '
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').Value = '3.76'
$ws.Range('E22').Value = '  +5.96%  '
$ws.Range('D23').Value = '73.42'
$ws.Range('E23').Value = '  -1.75%  '
$ws.Range('D24').Value = '261.74'
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '7.61'
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('D28').Value = '11.27'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('D30').Value = '22.45'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('D31').Value = '0.0958'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').Value = '37.42'
$ws.Range('E32').Value = '  -4.88%  '
$ws.Range('D33').Value = '167.59'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('D34').Value = '2.91'
$ws.Range('E34').Value = '  -1.32%  '
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = '4.74'
$ws.Range('E37').Value = '  -4.26%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '1.93'
$ws.Range('E38').Value = '  +12.07%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '4.04'
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('D40').Value = '2.98'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').Value = '0.0357'
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').Value = '97.96'
$ws.Range('E42').Value = '  -7.09%  '
$ws.Range('D43').Value = '70.79'
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.230'
$ws.Range('E44').Value = '  -4.10%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').Value = '12.97'
$ws.Range('E45').Value = '  -3.09%  '
$ws.Range('D46').Value = '1.861.57'
$ws.Range('E46').Value = '  +13.01%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '5.94'
$ws.Range('E48').Value = '  +2.77%  '
$ws.Range('D49').Value = '84.45'
$ws.Range('E49').Value = '  +6.99%  '
$ws.Range('D50').Value = '112.68'
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('D51').Value = '9.29'
$ws.Range('E51').Value = '  -0.52%  '

foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
